$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing row 4 (TC_EC_0003): Executed -> no
$ws.Range("D4").Value = "no"

# Insert a new row before current row 6 (END), shifting it down
$ws.Rows.Item(6).Insert()

# Fill in the new test case row
$ws.Range("A6").Value = "TC_EC_0005"
$ws.Range("B6").Value = "Mobile Registration"
$ws.Range("C6").Value = "Regression"
$ws.Range("D6").Value = "no"
